# Auto-generated edit script: apply Leviathan_Profits scheduled-runner price updates
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 40657
$ws.Range("J3").Value = 40657
$ws.Range("L3").Value = 40657
$ws.Range("N3").Value = -40885
$ws.Range("H38").Value = 707.1111
$ws.Range("I38").Value = 707.1111
$ws.Range("K38").Value = 2121.3333
$ws.Range("M38").Value = -1749.3333
$ws.Range("H41").Value = 947.8946999999999
$ws.Range("I41").Value = 1081.0714
$ws.Range("J41").Value = 575
$ws.Range("K41").Value = 1081.0714
$ws.Range("L41").Value = 575
$ws.Range("M41").Value = -641.0714
$ws.Range("N41").Value = -1455
$ws.Range("H62").Value = 9993.333000000001
$ws.Range("I62").Value = 5995
$ws.Range("K62").Value = 5995
$ws.Range("M62").Value = -5371
$ws.Range("H65").Value = 9993.333000000001
$ws.Range("I65").Value = 5995
$ws.Range("K65").Value = 29975
$ws.Range("M65").Value = -26855
$ws.Range("H80").Value = 1958.3043
$ws.Range("I80").Value = 995.25
$ws.Range("J80").Value = 3008.9092
$ws.Range("K80").Value = 2985.75
$ws.Range("L80").Value = 9026.7276
$ws.Range("M80").Value = -1987.75
$ws.Range("N80").Value = -11022.7276
$ws.Range("H83").Value = 1958.3043
$ws.Range("I83").Value = 995.25
$ws.Range("J83").Value = 3008.9092
$ws.Range("K83").Value = 8957.25
$ws.Range("L83").Value = 27080.1828
$ws.Range("M83").Value = -3965.25
$ws.Range("N83").Value = -37064.1828
$ws.Range("H97").Value = 1332
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 1332
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 3996
$ws.Range("M97").ClearContents()
$ws.Range("N97").Value = -4988
$ws.Range("H99").Value = 90924590
$ws.Range("I99").Value = 24053
$ws.Range("J99").Value = 250000540
$ws.Range("K99").Value = 72159
$ws.Range("L99").Value = 750001620
$ws.Range("M99").Value = -70661
$ws.Range("N99").Value = -750004616
$ws.Range("H101").Value = 1239.8462
$ws.Range("I101").Value = 1239.8462
$ws.Range("J101").Value = 0
$ws.Range("K101").Value = 3719.5386
$ws.Range("L101").Value = 0
$ws.Range("M101").Value = -2097.5386
$ws.Range("N101").ClearContents()
$ws.Range("H102").Value = 40657
$ws.Range("J102").Value = 40657
$ws.Range("L102").Value = 40657
$ws.Range("N102").Value = -47147
$ws.Range("H103").Value = 62500984
$ws.Range("J103").Value = 71429624
$ws.Range("L103").Value = 214288872
$ws.Range("N103").Value = -214290044
$ws.Range("H129").Value = 1814.826
$ws.Range("I129").Value = 712
$ws.Range("J129").Value = 3017.9092
$ws.Range("K129").Value = 2136
$ws.Range("L129").Value = 9053.7276
$ws.Range("M129").Value = 2864
$ws.Range("N129").Value = -19053.7276
$ws.Range("H132").Value = 1775.0938
$ws.Range("I132").Value = 1829.3572
$ws.Range("J132").Value = 1395.25
$ws.Range("K132").Value = 5488.071599999999
$ws.Range("L132").Value = 4185.75
$ws.Range("M132").Value = -2958.071599999999
$ws.Range("N132").Value = -9245.75
$ws.Range("H138").Value = 1798.1212
$ws.Range("I138").Value = 1124.5186
$ws.Range("J138").Value = 2264.4614
$ws.Range("K138").Value = 3373.5558
$ws.Range("L138").Value = 6793.3842
$ws.Range("M138").Value = 1766.4442
$ws.Range("N138").Value = -17073.3842

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6795.121
$ws.Range("I32").Value = 6379.3335
$ws.Range("J32").Value = 8666.166999999999
$ws.Range("K32").Value = 6379.3335
$ws.Range("L32").Value = 8666.166999999999
$ws.Range("M32").Value = -6092.3335
$ws.Range("N32").Value = -9240.166999999999
$ws.Range("H45").Value = 8962.111000000001
$ws.Range("I45").Value = 19167.834
$ws.Range("J45").Value = 3859.25
$ws.Range("K45").Value = 19167.834
$ws.Range("L45").Value = 3859.25
$ws.Range("M45").Value = -18790.834
$ws.Range("N45").Value = -4613.25
$ws.Range("H61").Value = 2360.7368
$ws.Range("I61").Value = 1304.8
$ws.Range("J61").Value = 3534
$ws.Range("K61").Value = 1304.8
$ws.Range("L61").Value = 3534
$ws.Range("M61").Value = -1092.8
$ws.Range("N61").Value = -3958
$ws.Range("H74").Value = 1770.1875
$ws.Range("I74").Value = 1770.1875
$ws.Range("K74").Value = 1770.1875
$ws.Range("M74").Value = -896.1875
$ws.Range("H77").Value = 1770.1875
$ws.Range("I77").Value = 1770.1875
$ws.Range("K77").Value = 8850.9375
$ws.Range("M77").Value = -4482.9375
$ws.Range("H102").Value = 4257.9165
$ws.Range("I102").Value = 4008.6365
$ws.Range("K102").Value = 4008.6365
$ws.Range("M102").Value = -2386.6365
$ws.Range("H110").Value = 1549.5
$ws.Range("I110").Value = 1549.5
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 1549.5
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = 495.5
$ws.Range("N110").ClearContents()
$ws.Range("H122").Value = 2341.3333
$ws.Range("I122").Value = 2341.3333
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 7023.999899999999
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -4573.999899999999
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value = 2345.8235
$ws.Range("I132").Value = 2201.8708
$ws.Range("K132").Value = 6605.6124
$ws.Range("M132").Value = -4075.6124
$ws.Range("H136").Value = 2360.7368
$ws.Range("I136").Value = 1304.8
$ws.Range("J136").Value = 3534
$ws.Range("K136").Value = 3914.4
$ws.Range("L136").Value = 10602
$ws.Range("M136").Value = -1364.4
$ws.Range("N136").Value = -15702

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 6311.2334
$ws.Range("I20").Value = 5989.9565
$ws.Range("K20").Value = 5989.9565
$ws.Range("M20").Value = -5742.9565
$ws.Range("H80").Value = 898.4286
$ws.Range("J80").Value = 1111.2222
$ws.Range("L80").Value = 1111.2222
$ws.Range("N80").Value = -3107.2222
$ws.Range("H83").Value = 898.4286
$ws.Range("J83").Value = 1111.2222
$ws.Range("L83").Value = 5556.111
$ws.Range("N83").Value = -15540.111
$ws.Range("H105").Value = 50010
$ws.Range("I105").Value = 50010
$ws.Range("K105").Value = 50010
$ws.Range("M105").Value = -48263
$ws.Range("H134").Value = 1317.871
$ws.Range("I134").Value = 1293
$ws.Range("K134").Value = 3879
$ws.Range("M134").Value = -1344

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 55559052
$ws.Range("I7").Value = 90911944
$ws.Range("J7").Value = 4514.2856
$ws.Range("K7").Value = 90911944
$ws.Range("L7").Value = 4514.2856
$ws.Range("M7").Value = -90911831
$ws.Range("N7").Value = -4740.2856
$ws.Range("H31").Value = 12787.5
$ws.Range("I31").Value = 1833.7084
$ws.Range("J31").Value = 78510.25
$ws.Range("K31").Value = 1833.7084
$ws.Range("L31").Value = 78510.25
$ws.Range("M31").Value = -1538.7084
$ws.Range("N31").Value = -79100.25
$ws.Range("H34").Value = 12787.5
$ws.Range("I34").Value = 1833.7084
$ws.Range("J34").Value = 78510.25
$ws.Range("K34").Value = 1833.7084
$ws.Range("L34").Value = 78510.25
$ws.Range("M34").Value = -1631.7084
$ws.Range("N34").Value = -78914.25
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()
$ws.Range("H58").Value = 2434.6667
$ws.Range("J58").Value = 2999
$ws.Range("L58").Value = 2999
$ws.Range("N58").Value = -3405
$ws.Range("H60").Value = 14000
$ws.Range("J60").Value = 14000
$ws.Range("L60").Value = 14000
$ws.Range("N60").Value = -15022
$ws.Range("H61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()
$ws.Range("H93").Value = 12458.111
$ws.Range("I93").Value = 8137.8335
$ws.Range("J93").Value = 21098.666
$ws.Range("K93").Value = 8137.8335
$ws.Range("L93").Value = 21098.666
$ws.Range("M93").Value = -6265.8335
$ws.Range("N93").Value = -24842.666
$ws.Range("H103").Value = 24341.334
$ws.Range("I103").Value = 21512
$ws.Range("J103").Value = 30000
$ws.Range("K103").Value = 21512
$ws.Range("L103").Value = 30000
$ws.Range("M103").Value = -20340
$ws.Range("N103").Value = -32344
$ws.Range("H106").Value = 64059.168
$ws.Range("J106").Value = 64059.168
$ws.Range("L106").Value = 64059.168
$ws.Range("N106").Value = -66583.16800000001
$ws.Range("H112").Value = 29992.5
$ws.Range("J112").Value = 29992.5
$ws.Range("L112").Value = 29992.5
$ws.Range("N112").Value = -32946.5
$ws.Range("H122").Value = 3960
$ws.Range("I122").Value = 3940
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 11820
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -9370
$ws.Range("N122").Value = -16900
$ws.Range("H134").Value = 1685.421
$ws.Range("I134").Value = 1236
$ws.Range("K134").Value = 3708
$ws.Range("M134").Value = -1173
$ws.Range("H136").Value = 2434.6667
$ws.Range("J136").Value = 2999
$ws.Range("L136").Value = 8997
$ws.Range("N136").Value = -14097
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()
$ws.Range("H141").Value = 264611.12
$ws.Range("J141").Value = 264611.12
$ws.Range("L141").Value = 264611.12
$ws.Range("N141").Value = -274971.12

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1877164
$ws.Range("I4").Value = 1000408.4
$ws.Range("J4").Value = 4069053
$ws.Range("K4").Value = 3001225.2
$ws.Range("L4").Value = 12207159
$ws.Range("M4").Value = -3001113.2
$ws.Range("N4").Value = -12207383
$ws.Range("H12").Value = 140.08333
$ws.Range("J12").Value = 152
$ws.Range("L12").Value = 456
$ws.Range("N12").Value = -802
$ws.Range("H34").Value = 1879.3529
$ws.Range("I34").Value = 219.83333
$ws.Range("J34").Value = 2784.5454
$ws.Range("K34").Value = 659.49999
$ws.Range("L34").Value = 8353.636200000001
$ws.Range("M34").Value = -575.49999
$ws.Range("N34").Value = -8521.636200000001
$ws.Range("H38").Value = 31.25
$ws.Range("I38").Value = 31.25
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 93.75
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = 253.25
$ws.Range("N38").ClearContents()
$ws.Range("H51").Value = 2599.5
$ws.Range("I51").Value = 2599.5
$ws.Range("K51").Value = 7798.5
$ws.Range("M51").Value = -7338.5
$ws.Range("H60").Value = 808.8461
$ws.Range("I60").Value = 459.85715
$ws.Range("J60").Value = 1216
$ws.Range("K60").Value = 1379.57145
$ws.Range("L60").Value = 3648
$ws.Range("M60").Value = -1128.57145
$ws.Range("N60").Value = -4150
$ws.Range("H63").Value = 1466
$ws.Range("J63").Value = 899
$ws.Range("L63").Value = 2697
$ws.Range("N63").Value = -4195
$ws.Range("H66").Value = 1466
$ws.Range("J66").Value = 899
$ws.Range("L66").Value = 8091
$ws.Range("N66").Value = -15579
$ws.Range("H104").Value = 3753.6365
$ws.Range("J104").Value = 4165.5557
$ws.Range("L104").Value = 12496.6671
$ws.Range("N104").Value = -17738.6671
$ws.Range("H107").Value = 670.2308
$ws.Range("J107").Value = 668.44446
$ws.Range("L107").Value = 2005.33338
$ws.Range("N107").Value = -5845.33338
$ws.Range("H113").Value = 1486.6154
$ws.Range("J113").Value = 1602.5
$ws.Range("L113").Value = 4807.5
$ws.Range("N113").Value = -9147.5
$ws.Range("H122").Value = 1197
$ws.Range("I122").Value = 715
$ws.Range("J122").Value = 1920
$ws.Range("K122").Value = 6435
$ws.Range("L122").Value = 17280
$ws.Range("M122").Value = -3985
$ws.Range("N122").Value = -22180
$ws.Range("H134").Value = 19870
$ws.Range("I134").Value = 3971.6667
$ws.Range("K134").Value = 11915.0001
$ws.Range("M134").Value = -6845.000100000001
$ws.Range("H140").Value = 2445.3635
$ws.Range("I140").Value = 2389.9
$ws.Range("J140").Value = 3000
$ws.Range("K140").Value = 7169.700000000001
$ws.Range("L140").Value = 9000
$ws.Range("M140").Value = -1989.700000000001
$ws.Range("N140").Value = -19360

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("N53").ClearContents()
$ws.Range("H70").Value = 6829
$ws.Range("I70").Value = 6829
$ws.Range("K70").Value = 6829
$ws.Range("M70").Value = -6559
$ws.Range("H73").Value = 6829
$ws.Range("I73").Value = 6829
$ws.Range("K73").Value = 6829
$ws.Range("M73").Value = -5893
$ws.Range("H102").Value = 2643.12
$ws.Range("I102").Value = 2743.1738
$ws.Range("K102").Value = 2743.1738
$ws.Range("M102").Value = -1121.1738
$ws.Range("H113").Value = 4308.4443
$ws.Range("I113").Value = 4158.8
$ws.Range("J113").Value = 4495.5
$ws.Range("K113").Value = 4158.8
$ws.Range("L113").Value = 4495.5
$ws.Range("M113").Value = -1988.8
$ws.Range("N113").Value = -8835.5
$ws.Range("H132").Value = 1886.5333
$ws.Range("I132").Value = 1877.625
$ws.Range("K132").Value = 5632.875
$ws.Range("M132").Value = -3102.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1781.6666
$ws.Range("I22").Value = 1781.6666
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 1781.6666
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -1486.6666
$ws.Range("N22").ClearContents()
$ws.Range("H27").Value = 1781.6666
$ws.Range("I27").Value = 1781.6666
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 1781.6666
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -1674.6666
$ws.Range("N27").ClearContents()
$ws.Range("H46").Value = 62359.715
$ws.Range("I46").Value = 212361.5
$ws.Range("K46").Value = 212361.5
$ws.Range("M46").Value = -212173.5
$ws.Range("H61").Value = 88962.69
$ws.Range("I61").Value = 60213.234
$ws.Range("K61").Value = 60213.234
$ws.Range("M61").Value = -60011.234
$ws.Range("H98").Value = 30327.5
$ws.Range("J98").Value = 30327.5
$ws.Range("L98").Value = 30327.5
$ws.Range("N98").Value = -36317.5
$ws.Range("H113").Value = 88962.69
$ws.Range("I113").Value = 60213.234
$ws.Range("K113").Value = 60213.234
$ws.Range("M113").Value = -58043.234
$ws.Range("H122").Value = 115069.72
$ws.Range("I122").Value = 156891.47
$ws.Range("J122").Value = 6333.2
$ws.Range("K122").Value = 470674.41
$ws.Range("L122").Value = 18999.6
$ws.Range("M122").Value = -468224.41
$ws.Range("N122").Value = -23899.6
$ws.Range("H127").Value = 57602
$ws.Range("I127").Value = 24999
$ws.Range("J127").Value = 64122.6
$ws.Range("K127").Value = 24999
$ws.Range("L127").Value = 64122.6
$ws.Range("M127").Value = -20039
$ws.Range("N127").Value = -74042.60000000001
$ws.Range("H132").Value = 5330.5264
$ws.Range("I132").Value = 5098.923
$ws.Range("K132").Value = 15296.769
$ws.Range("M132").Value = -12766.769
$ws.Range("H135").Value = 116840.5
$ws.Range("J135").Value = 116840.5
$ws.Range("L135").Value = 116840.5
$ws.Range("N135").Value = -126980.5
$ws.Range("H136").Value = 5231.6665
$ws.Range("I136").Value = 4465
$ws.Range("K136").Value = 13395
$ws.Range("M136").Value = -10845

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 42249.75
$ws.Range("H81").Value = 3147.6667
$ws.Range("I81").Value = 3147.6667
$ws.Range("K81").Value = 6295.3334
$ws.Range("M81").Value = -5234.3334
$ws.Range("H84").Value = 3147.6667
$ws.Range("I84").Value = 3147.6667
$ws.Range("K84").Value = 31476.667
$ws.Range("M84").Value = -26172.667
$ws.Range("H113").Value = 213.2258
$ws.Range("J113").Value = 250.33333
$ws.Range("L113").Value = 750.99999
$ws.Range("N113").Value = -5090.99999
$ws.Range("H122").Value = 1791.625
$ws.Range("I122").Value = 1690.4286
$ws.Range("K122").Value = 5071.2858
$ws.Range("M122").Value = -2621.2858
$ws.Range("H132").Value = 2683.1482
$ws.Range("I132").Value = 2902.45
$ws.Range("J132").Value = 2056.5715
$ws.Range("K132").Value = 8707.349999999999
$ws.Range("L132").Value = 6169.7145
$ws.Range("M132").Value = -6177.349999999999
$ws.Range("N132").Value = -11229.7145
$ws.Range("H136").Value = 7563.8423
$ws.Range("I136").Value = 7428.5
$ws.Range("K136").Value = 22285.5
$ws.Range("M136").Value = -19735.5

Write-Host "Applied 421 value updates and 11 clears across 8 sheets"